# Refresh the coin-ranking snapshot: updated prices / 1h volume %,
# the coin roster shift (rows 6-17), and the "Hora" snapshot hour (15 -> 16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'302.37"
$ws.Range("E2").Value = "'2.00%"
$ws.Range("G2").Value = "'16"
# Row 3
$ws.Range("D3").Value = "'43.20"
$ws.Range("E3").Value = "'4.58%"
$ws.Range("G3").Value = "'16"
# Row 4
$ws.Range("D4").Value = "'5.075"
$ws.Range("E4").Value = "'0.41%"
$ws.Range("G4").Value = "'16"
# Row 5
$ws.Range("D5").Value = "'0.07689"
$ws.Range("E5").Value = "'2.84%"
$ws.Range("G5").Value = "'16"
# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.629"
$ws.Range("E6").Value = "'3.42%"
$ws.Range("G6").Value = "'16"
# Row 7
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'1.055"
$ws.Range("E7").Value = "'13.57%"
$ws.Range("G7").Value = "'16"
# Row 8
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "'0.1267"
$ws.Range("E8").Value = "'5.36%"
$ws.Range("G8").Value = "'16"
# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1860"
$ws.Range("E9").Value = "'1.85%"
$ws.Range("G9").Value = "'16"
# Row 10
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.09141"
$ws.Range("E10").Value = "'3.84%"
$ws.Range("G10").Value = "'16"
# Row 11
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04164"
$ws.Range("E11").Value = "'-2.84%"
$ws.Range("G11").Value = "'16"
# Row 12
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.1048"
$ws.Range("E12").Value = "'-0.26%"
$ws.Range("G12").Value = "'16"
# Row 13
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001268"
$ws.Range("E13").Value = "'-1.30%"
$ws.Range("G13").Value = "'16"
# Row 14
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.005743"
$ws.Range("E14").Value = "'-4.53%"
$ws.Range("G14").Value = "'16"
# Row 15
$ws.Range("B15").Value = "UpBots"
$ws.Range("C15").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D15").Value = "'0.007430"
$ws.Range("E15").Value = "'1,897.56%"
$ws.Range("G15").Value = "'16"
# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.347"
$ws.Range("E16").Value = "'-0.35%"
$ws.Range("G16").Value = "'16"
# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.424"
$ws.Range("E17").Value = "'1.16%"
$ws.Range("G17").Value = "'16"
# Row 18
$ws.Range("E18").Value = "'-1.96%"
$ws.Range("G18").Value = "'16"
# Row 19
$ws.Range("D19").Value = "'0.3356"
$ws.Range("E19").Value = "'1.62%"
$ws.Range("G19").Value = "'16"
# Row 20
$ws.Range("D20").Value = "'8.645"
$ws.Range("E20").Value = "'7.36%"
$ws.Range("G20").Value = "'16"
# Row 21
$ws.Range("D21").Value = "'0.1367"
$ws.Range("E21").Value = "'-0.87%"
$ws.Range("G21").Value = "'16"
# Row 22
$ws.Range("D22").Value = "'0.3170"
$ws.Range("E22").Value = "'6.94%"
$ws.Range("G22").Value = "'16"
# Row 23
$ws.Range("D23").Value = "'0.04182"
$ws.Range("E23").Value = "'3.84%"
$ws.Range("G23").Value = "'16"
# Row 24
$ws.Range("D24").Value = "'0.001281"
$ws.Range("E24").Value = "'1.26%"
$ws.Range("G24").Value = "'16"
# Row 25
$ws.Range("D25").Value = "'0.004432"
$ws.Range("E25").Value = "'14.65%"
$ws.Range("G25").Value = "'16"
# Row 26
$ws.Range("E26").Value = "'9.63%"
$ws.Range("G26").Value = "'16"
# Row 27
$ws.Range("G27").Value = "'16"
# Row 28
$ws.Range("G28").Value = "'16"
# Row 29
$ws.Range("G29").Value = "'16"
# Row 30
$ws.Range("G30").Value = "'16"
# Row 31
$ws.Range("G31").Value = "'16"
# Row 32
$ws.Range("G32").Value = "'16"
# Row 33
$ws.Range("G33").Value = "'16"
# Row 34
$ws.Range("G34").Value = "'16"
# Row 35
$ws.Range("G35").Value = "'16"
# Row 36
$ws.Range("G36").Value = "'16"
# Row 37
$ws.Range("G37").Value = "'16"
# Row 38
$ws.Range("D38").Value = "'0.02484"
$ws.Range("E38").Value = "'3.33%"
$ws.Range("G38").Value = "'16"
# Row 39
$ws.Range("D39").Value = "'0.05293"
$ws.Range("E39").Value = "'2.19%"
$ws.Range("G39").Value = "'16"
# Row 40
$ws.Range("D40").Value = "'0.005932"
$ws.Range("E40").Value = "'-11.34%"
$ws.Range("G40").Value = "'16"
# Row 41
$ws.Range("D41").Value = "'0.007649"
$ws.Range("E41").Value = "'-1.61%"
$ws.Range("G41").Value = "'16"
# Row 42
$ws.Range("D42").Value = "'0.1351"
$ws.Range("E42").Value = "'2.41%"
$ws.Range("G42").Value = "'16"
# Row 43
$ws.Range("D43").Value = "'0.007365"
$ws.Range("G43").Value = "'16"
# Row 44
$ws.Range("D44").Value = "'0.007527"
$ws.Range("E44").Value = "'-3.85%"
$ws.Range("G44").Value = "'16"
# Row 45
$ws.Range("D45").Value = "'0.3009"
$ws.Range("E45").Value = "'-6.28%"
$ws.Range("G45").Value = "'16"
# Row 46
$ws.Range("D46").Value = "'0.00006717"
$ws.Range("E46").Value = "'6.11%"
$ws.Range("G46").Value = "'16"
# Row 47
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("G47").Value = "'16"
# Row 48
$ws.Range("D48").Value = "'0.04479"
$ws.Range("E48").Value = "'-3.50%"
$ws.Range("G48").Value = "'16"
# Row 49
$ws.Range("E49").Value = "'0.08%"
$ws.Range("G49").Value = "'16"
# Row 50
$ws.Range("D50").Value = "'0.00002095"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("G50").Value = "'16"
# Row 51
$ws.Range("D51").Value = "'0.0001995"
$ws.Range("E51").Value = "'-0.20%"
$ws.Range("G51").Value = "'16"
